$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (row 11)
$ws.Range("R11").Value = 0
$ws.Range("T11").Value = 0

# Row 21
$ws.Range("T21").Value = 0

# Row 23
$ws.Range("R23").Value = 0
$ws.Range("T23").Value = 5

# Row 26
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 0
$ws.Range("N26").Value = 2
$ws.Range("O26").Value = 5
$ws.Range("P26").Value = 0
$ws.Range("R26").Value = 0
$ws.Range("S26").Value = 0
$ws.Range("T26").Value = 5

# Update selection/view
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("J26").Select()
